$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 533 - this shifts the existing rows
# 533-601 down to 534-602 (dimension grows from A1:R601 to A1:R602).
$ws.Rows.Item(533).Insert()

# Populate the newly inserted row 533 with the new weekly price record.
$ws.Range("A533").Value = 4
$ws.Range("B533").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C533").Value = "Los Lagos"
$ws.Range("D533").Value = 45124
$ws.Range("E533").Value = 10
$ws.Range("F533").Value = 100114013
$ws.Range("G533").Value = "Zanahoria"
$ws.Range("H533").Value = "Sin especificar"
$ws.Range("I533").Value = "Primera"
$ws.Range("J533").Value = 400
$ws.Range("K533").Value = 8000
$ws.Range("L533").Value = 8000
$ws.Range("M533").Value = 8000
$ws.Range("N533").Value = "`$/saco 20 kilos"
$ws.Range("O533").Value = "Provincia de Llanquihue"
$ws.Range("P533").Value = 400
$ws.Range("Q533").Value = 20
$ws.Range("R533").Value = "Hortaliza"
